$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------
# New "Alternativni tok/korak" blocks (rows 21-60).
#
# Cell values are written in a deliberate order: every brand-new piece of
# text is assigned to a cell for the first time in the same order the
# original author must have typed it in (so the shared-strings table the
# engine builds ends up in the same order as the target workbook), and only
# afterwards do we fill in the cells that reuse already-existing strings
# ("Neregistrovani korisnik", "Registracij korisnika(sistem)", the password
# validation text, etc).
# ----------------------------------------------------------------------------
$ws.Range("A21").Value = "Alternativni tok 2:"
$ws.Range("A22").Value = "Na koraku 8 nije unesena validna lozinka, korisnik odustaje od registracije"
$ws.Range("A29").Value = "Alternativni tok 3:"
$ws.Range("B32").Value = "1. Nije uneseno validno ime (sadrži nevalidne karaktere)"
$ws.Range("B34").Value = "3. Korisniku se nudi mogućnost `nponovnog unosa "
$ws.Range("B42").Value = "3. Korisniku se nudi mogućnost `nponovnog unosa "
$ws.Range("B51").Value = "3. Korisniku se nudi mogućnost `nponovnog unosa "
$ws.Range("B59").Value = "3. Korisniku se nudi mogućnost `nponovnog unosa "
$ws.Range("A35").Value = "4. Nastavak na koraku 2"
$ws.Range("A37").Value = "Alternativni korak 4:"
$ws.Range("B40").Value = "1. Nije uneseno validno korisničko ime (sadrži nevalidne karaktere)"
$ws.Range("A38").Value = "Na koraku 6  nije uneseno validno korisničko ime"
$ws.Range("A43").Value = "4. Nastavak na koraku 6"
$ws.Range("A30").Value = "Na koraku 2 nije uneseno validno ime (analogno i za prezime, s tim da se vraća na korak 4.)"
$ws.Range("A46").Value = "Alternativni korak 5:"
$ws.Range("A47").Value = "Na koraku 10 nije unesen validan datum rođenja"
$ws.Range("B49").Value = "1. Nije unesen validan datum"
$ws.Range("A52").Value = "4. Nastavak na koraku 10"
$ws.Range("A54").Value = "Alternativni korak 6:"
$ws.Range("A55").Value = "Na koraku 12 nije unesen validan email"
$ws.Range("B57").Value = "1. Nije unesen validan email"
$ws.Range("A60").Value = "4. Nastavak na koraku 12"

# Header rows of the five new mini-tables (reuse the existing header strings).
$ws.Range("B24").Value = "Registracij korisnika(sistem)"
$ws.Range("B31").Value = "Registracij korisnika(sistem)"
$ws.Range("B39").Value = "Registracij korisnika(sistem)"
$ws.Range("B48").Value = "Registracij korisnika(sistem)"
$ws.Range("B56").Value = "Registracij korisnika(sistem)"
$ws.Range("A24").Value = "Neregistrovani korisnik"
$ws.Range("A31").Value = "Neregistrovani korisnik"
$ws.Range("A39").Value = "Neregistrovani korisnik"
$ws.Range("A48").Value = "Neregistrovani korisnik"
$ws.Range("A56").Value = "Neregistrovani korisnik"

# Remaining cells that reuse already-existing shared strings.
$ws.Range("B25").Value = "1. Nije unesena validna lozinka `n(dozvoljena slova i brojevi, pri čemu barem jedno slovo mora biti veliko slovo te lozinka mora sadržati barem jednu cifru)"
$ws.Range("B26").Value = "2. Obavještavanje korisnika o problemu"
$ws.Range("B33").Value = "2. Obavještavanje korisnika o problemu"
$ws.Range("B41").Value = "2. Obavještavanje korisnika o problemu"
$ws.Range("B50").Value = "2. Obavještavanje korisnika o problemu"
$ws.Range("B58").Value = "2. Obavještavanje korisnika o problemu"
$ws.Range("B27").Value = "3. Korisniku se nudi mogućnost `nponovnog unosa lozinke"

# ----------------------------------------------------------------------------
# Wrap text on the long explanatory cells and give their rows the same
# height used by the equivalent cells in the first "Alternativni tok" block.
# ----------------------------------------------------------------------------
$ws.Range("B25").WrapText = $true
$ws.Rows(25).RowHeight = 75

$ws.Range("B27").WrapText = $true
$ws.Rows(27).RowHeight = 30

$ws.Range("B32").WrapText = $true
$ws.Rows(32).RowHeight = 30

$ws.Range("B34").WrapText = $true
$ws.Rows(34).RowHeight = 30

$ws.Range("B40").WrapText = $true
$ws.Rows(40).RowHeight = 30

$ws.Range("B42").WrapText = $true
$ws.Rows(42).RowHeight = 30

$ws.Range("B49").WrapText = $true

$ws.Range("B51").WrapText = $true
$ws.Rows(51).RowHeight = 30

$ws.Range("B57").WrapText = $true

$ws.Range("B59").WrapText = $true
$ws.Rows(59).RowHeight = 30

# ----------------------------------------------------------------------------
# Turn the five new two-column blocks into tables, just like the existing
# "Tabelle3" table at A15:B19.
# ----------------------------------------------------------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A24:B27"), 0, 1)
$lo.Name = "Tabelle33"

$lo = $ws.ListObjects.Add(1, $ws.Range("A31:B35"), 0, 1)
$lo.Name = "Tabelle335"

$lo = $ws.ListObjects.Add(1, $ws.Range("A39:B43"), 0, 1)
$lo.Name = "Tabelle5"

$lo = $ws.ListObjects.Add(1, $ws.Range("A48:B52"), 0, 1)
$lo.Name = "Tabelle58"

$lo = $ws.ListObjects.Add(1, $ws.Range("A56:B60"), 0, 1)
$lo.Name = "Tabelle589"

# ----------------------------------------------------------------------------
# Match the final selection/scroll state recorded in the target workbook.
# ----------------------------------------------------------------------------
$ws.Range("A60").Select()
$excel.ActiveWindow.ScrollRow = 35
